# Apply cell value updates from the coinranking data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.920.06"
$ws.Range("E2").Value = "  +2.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.111.79"
$ws.Range("E3").Value = "  +5.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.97"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.14"
$ws.Range("E6").Value = "  +6.96%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.105.20"
$ws.Range("E8").Value = "  +5.36%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("E12").Value = "  +5.19%  "
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.32"
$ws.Range("E15").Value = "  -0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.627.90"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.889.85"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.110.13"
$ws.Range("E19").Value = "  +5.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.14"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "484.48"
$ws.Range("E21").Value = "  +8.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("E22").Value = "  +3.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.53"
$ws.Range("E23").Value = "  +3.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.17"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +5.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.17"
$ws.Range("E26").Value = "  +7.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.40"
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.98"
$ws.Range("E32").Value = "  +6.50%  "
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +3.30%  "
$ws.Range("E37").Value = "  +3.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.01"
$ws.Range("E38").Value = "  +6.25%  "
$ws.Range("E39").Value = "  +7.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.21"
$ws.Range("E40").Value = "  +2.13%  "
$ws.Range("E41").Value = "  +4.90%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.81"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0361"
$ws.Range("E45").Value = "  +2.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.838.96"
$ws.Range("E46").Value = "  +5.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "382.60"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "135.48"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.05"
$ws.Range("E50").Value = "  +4.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("E51").Value = "  +2.80%  "
